# Correct the CAS registry number typo for Hemoglobin ("9908-02-0" -> "9008-02-0")
# and keep the CAS-sorted data in order (the corrected value now sorts before
# "92-31-9", so the Hemoglobin row moves up one position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldCas = '"9908-02-0"'
$newCas = '"9008-02-0"'

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Locate the row holding the typo'd CAS number in column A.
$target = $ws.Range("A1:A" + $lastRow).Find($oldCas)

if ($target -ne $null) {
    $row = $target.Row

    # Fix the typo in place first.
    $ws.Cells.Item($row, 1).Value = $newCas

    # Walk the corrected value up past any preceding rows it now sorts before
    # (text/alphabetic comparison, matching the existing CAS-sorted column A),
    # swapping whole rows so the rest of the sheet stays untouched.
    while ($row -gt 2) {
        $prevCas = $ws.Cells.Item($row - 1, 1).Value2

        if (-not ($newCas -clt $prevCas)) {
            break
        }

        $thisRowVals = @()
        $prevRowVals = @()
        for ($c = 1; $c -le 6; $c++) {
            $thisRowVals += ,$ws.Cells.Item($row, $c).Value2
            $prevRowVals += ,$ws.Cells.Item($row - 1, $c).Value2
        }

        for ($c = 1; $c -le 6; $c++) {
            $ws.Cells.Item($row - 1, $c).Value = $thisRowVals[$c - 1]
            $ws.Cells.Item($row, $c).Value = $prevRowVals[$c - 1]
        }

        $row = $row - 1
    }

    $ws.Cells.Item($lastRow, 1).Select() | Out-Null
}
